# Mobile Web testing for Chrome and Safari
# Applies updates to the "ECS" worksheet: corrects/reorders the
# LogIn Android / LogIn IOS rows (10-11) and appends two new rows
# (12-13) for Mobile Chrome and Mobile Safari testing.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ECS")

# --- Row 10: LogIn Android tests -----------------------------------
$ws.Range("A10").Value = "Regression, Sanity"
$ws.Range("B10").Value = "LogIn Android tests"
$ws.Range("C10").Value = "tests.mobile.ECSScriptsAndroid"
$ws.Range("D10").Value = "'009"
$ws.Range("E10").Value = "TC01_Android_logInToECSWithValidID"
$ws.Range("F10").Value = "env,browser"
$ws.Range("G10").Value = "ECS_TEST,Android"
$ws.Range("H10").Value = "N"
$ws.Range("I10").Value = "null"

# --- Row 11: LogIn IOS tests ----------------------------------------
$ws.Range("A11").Value = "Regression, Sanity"
$ws.Range("B11").Value = "LogIn IOS tests"
$ws.Range("C11").Value = "tests.mobile.ECSScriptsIOS"
$ws.Range("D11").Value = "'010"
$ws.Range("E11").Value = "TC01_IOS_logInToECSWithValidID"
$ws.Range("F11").Value = "env,browser"
$ws.Range("G11").Value = "ECS_TEST,IOS"
$ws.Range("H11").Value = "N"
$ws.Range("I11").Value = "null"

# --- Row 12: LogIn Mobile Chrome Testing -----------------------------
$ws.Range("A12").Value = "Regression, Sanity"
$ws.Range("B12").Value = "LogIn Mobile Chrome Testing"
$ws.Range("C12").Value = "tests.web.ECSScripts"
$ws.Range("D12").Value = "'011"
$ws.Range("E12").Value = "TC01_logInToECSWithValidID"
$ws.Range("F12").Value = "env,browser"
$ws.Range("G12").Value = "ECS_TEST,android_chrome"
$ws.Range("H12").Value = "N"
$ws.Range("I12").Value = "null"

# --- Row 13: LogIn Mobile Safari Testing -----------------------------
$ws.Range("A13").Value = "Regression, Sanity"
$ws.Range("B13").Value = "LogIn Mobile Safari Testing"
$ws.Range("C13").Value = "tests.web.ECSScripts"
$ws.Range("D13").Value = "'012"
$ws.Range("E13").Value = "TC01_logInToECSWithValidID"
$ws.Range("F13").Value = "env,browser"
$ws.Range("G13").Value = "ECS_TEST,ios_safari"
$ws.Range("H13").Value = "Y"
$ws.Range("I13").Value = "null"

# --- column G width (widened to fit new content) ---------------------
$ws.Columns.Item(7).ColumnWidth = 22.17

# --- view state: scroll position + selection --------------------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 3
$excel.ActiveWindow.ScrollColumn = 2
$ws.Range("F18").Select()
